$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text updates: "In Translation" -> "Ready for handoff" ---------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamp updates -----------------------------------------------------
# Overview!G2 and de-de!H2 shared the same original string
# ("2016-08-31 17:08:50") which becomes "2016-08-31 17:09:28".
$wsOverview.Range("G2").Value = "2016-08-31 17:09:28"
$wsDeDe.Range("H2").Value     = "2016-08-31 17:09:28"

# zh-cn!H2 had its own original string ("2016-08-31 17:08:45") which
# becomes "2016-08-31 17:09:23".
$wsZhCn.Range("H2").Value = "2016-08-31 17:09:23"

# --- Column width updates ---------------------------------------------------
# The "Status" columns widen (13.4101845877511 -> 17.2159881591797 chars)
# to fit the longer "Ready for handoff" text. The host engine quantizes
# ColumnWidth to steps of 1/6 character, so 16.3333... is the input that
# lands closest to the target stored width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.333333333333332
